# Simplify the solution and add comments
#
# The sample sheet used to show a "PrevTime" RTD topic (column E) plus a
# helper column F that diffed "local time" against "PrevTime" to prove the
# value was refreshing. That's no longer needed, so this script:
#   1. Drops the PrevTime column (E) entirely.
#   2. Clears the now-orphaned diff formula in column F, but keeps the cell
#      (and its h:mm:ss style) around as a placeholder.
#   3. Un-hides row 6 and resets the AutoFilter so it only covers the
#      remaining A:D columns, with no filter criteria applied.
#   4. Keeps the workbook's _FilterDatabase defined name in sync with the
#      smaller AutoFilter range.
#   5. Leaves the selection on B2 (top of the shrunk table) instead of the
#      old scrolled-down D8 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1 & 3: clear out any active filter criteria and show every row again ---
$ws.AutoFilterMode = $false
$ws.Rows.Item(6).Hidden = $false

# --- 1: remove the "PrevTime" header/data column (E) ---
$ws.Range("E1").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E6:E10").ClearContents()

# --- 2: clear the leftover "local time - PrevTime" helper formulas in F,
#        but leave the cells (and their time-format style) in place ---
$ws.Range("F1").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("F6:F10").ClearContents()

# --- 3: re-apply a plain AutoFilter over the smaller A5:D10 range ---
$null = $ws.Range("A5:D10").AutoFilter()

# --- 4: keep the _FilterDatabase defined name pointed at the new range ---
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=" + $ws.Name + "!" + '$A$5:$D$10'
    }
}

# --- 5: move the selection back to the top of the table ---
$ws.Range("B2").Select()
